$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# --- Update raw numbers table (rows 6-8) with new SA results ---
$ws.Range("B6").Value = 54.5
$ws.Range("C6").Value = 55
$ws.Range("D6").Value = 60
$ws.Range("E6").Value = 45.5
$ws.Range("F6").Value = 45.5
$ws.Range("G6").Value = 61.5

$ws.Range("B7").Value = 55
$ws.Range("C7").Value = 53.5
$ws.Range("D7").Value = 58
$ws.Range("E7").Value = 51.5
$ws.Range("F7").Value = 49
$ws.Range("G7").Value = 58

$ws.Range("D8").Value = 58
$ws.Range("F8").Value = 63.5
$ws.Range("G8").Value = 53

# --- Insert a new "First Ipsos ->" row into the raw numbers table ---
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "First Ipsos ->"
$ws.Range("B14").Value = 55
$ws.Range("C14").Value = 52
$ws.Range("D14").Value = 58
$ws.Range("E14").Value = 50
$ws.Range("F14").Value = 59
$ws.Range("G14").Value = 55

# --- Insert a new row into the swing deviations table for the new poll ---
$ws.Rows.Item(27).Insert()
$ws.Range("A26").Value = "First Ipsos ->"
$ws.Range("C26").Formula = '=(C14-C$2)-($B14-$B$2)'
$ws.Range("D26").Formula = '=(D14-D$2)-($B14-$B$2)'
$ws.Range("E26").Formula = '=(E14-E$2)-($B14-$B$2)'
$ws.Range("F26").Formula = '=(F14-F$2)-($B14-$B$2)'
$ws.Range("G26").Formula = '=(G14-G$2)-($B14-$B$2)'
$ws.Range("H25").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("H26").Formula = '=F26*0.439183+G26*0.336042'

# --- Extend the weighted-average formula (row 28) to include the new poll ---
$ws.Range("C28").Formula = '=AVERAGE(C17,C17,AVERAGE(C18,C19,C20),AVERAGE(C21,C22,C23),AVERAGE(C24,C25),AVERAGE(C26,C17,C17,AVERAGE(C18,C19,C20),AVERAGE(C21,C22,C23),AVERAGE(C24,C25),C17,C17,AVERAGE(C18,C19,C20),AVERAGE(C21,C22,C23),AVERAGE(C24,C25)))'
$ws.Range("D28").Formula = '=AVERAGE(D17,D17,AVERAGE(D18,D19,D20),AVERAGE(D21,D22,D23),AVERAGE(D24,D25),AVERAGE(D26,D17,D17,AVERAGE(D18,D19,D20),AVERAGE(D21,D22,D23),AVERAGE(D24,D25),D17,D17,AVERAGE(D18,D19,D20),AVERAGE(D21,D22,D23),AVERAGE(D24,D25)))'
$ws.Range("E28").Formula = '=AVERAGE(E17,E17,AVERAGE(E18,E19,E20),AVERAGE(E21,E22,E23),AVERAGE(E24,E25),AVERAGE(E26,E17,E17,AVERAGE(E18,E19,E20),AVERAGE(E21,E22,E23),AVERAGE(E24,E25),E17,E17,AVERAGE(E18,E19,E20),AVERAGE(E21,E22,E23),AVERAGE(E24,E25)))'
$ws.Range("F28").Formula = '=AVERAGE(F17,F17,AVERAGE(F18,F19,F20),AVERAGE(F21,F22,F23),AVERAGE(F24,F25),AVERAGE(F26,F17,F17,AVERAGE(F18,F19,F20),AVERAGE(F21,F22,F23),AVERAGE(F24,F25),F17,F17,AVERAGE(F18,F19,F20),AVERAGE(F21,F22,F23),AVERAGE(F24,F25)))'
$ws.Range("G28").Formula = '=AVERAGE(G17,G17,AVERAGE(G18,G19,G20),AVERAGE(G21,G22,G23),AVERAGE(G24,G25),AVERAGE(G26,G17,G17,AVERAGE(G18,G19,G20),AVERAGE(G21,G22,G23),AVERAGE(G24,G25),G17,G17,AVERAGE(G18,G19,G20),AVERAGE(G21,G22,G23),AVERAGE(G24,G25)))'

# --- Update selection to match the author's final cursor position ---
$ws.Range("K26").Select()
